# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect a failed
# handback transform:
#   - "Overview" sheet: the status text "Ready for handoff" is replaced
#     (wherever that shared string is used) with "Handback transform failed".
#   - "zh-cn" and "de-de" sheets: the "Error Detail" column (P) gets a
#     widened column and a new error message filled into the last data row,
#     and that column is widened to fit the longer content.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: update status text ---
# E3/F3 (and any other cell sharing this string) hold "Ready for handoff";
# change the text to "Handback transform failed".
$ws_overview.Range("E3").Value = "Handback transform failed"
$ws_overview.Range("F3").Value = "Handback transform failed"

# The "Status" column (C) on the zh-cn/de-de sheets references the same
# shared string text, so it is updated to match automatically as well.
$ws_zhcn.Range("C3").Value = "Handback transform failed"
$ws_dede.Range("C3").Value = "Handback transform failed"

# --- zh-cn sheet: Error Detail column ---
$ws_zhcn.Range("P3").Value = "Handback file name: xkvz33xm.ikf is different with handoff file name: e74a2ec4-f7be-4d4b-a20f-84488dbaf2b3.57427849ddde75bba697d604749aa12bb5950c31.zh-cn."
# The saved OOXML column width is ColumnWidth + 0.8333333333333334, so subtract
# that offset here in order to land exactly on a stored width of 40.
$ws_zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet: Error Detail column ---
$ws_dede.Range("P3").Value = "Handback file name: xkvz33xm.ikf is different with handoff file name: e74a2ec4-f7be-4d4b-a20f-84488dbaf2b3.57427849ddde75bba697d604749aa12bb5950c31.de-de."
$ws_dede.Columns.Item(16).ColumnWidth = 39.166666666666664
